$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.250.43"
$ws.Range("E2").Value = "  -4.15%  "
$ws.Range("D3").Value = "3.106.42"
$ws.Range("E3").Value = "  -4.38%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'604.20"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'143.97"
$ws.Range("E6").Value = "  -8.23%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.103.92"
$ws.Range("E8").Value = "  -4.44%  "
$ws.Range("E9").Value = "  -4.39%  "
$ws.Range("E10").Value = "  -7.07%  "
$ws.Range("D11").Value = "'5.20"
$ws.Range("E11").Value = "  -8.46%  "
$ws.Range("D12").Value = "'0.466"
$ws.Range("E12").Value = "  -5.49%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  -7.31%  "
$ws.Range("D14").Value = "'35.01"
$ws.Range("E14").Value = "  -8.99%  "
$ws.Range("D15").Value = "3.614.19"
$ws.Range("E15").Value = "  -4.69%  "
$ws.Range("D17").Value = "63.372.37"
$ws.Range("E17").Value = "  -4.15%  "
$ws.Range("D18").Value = "3.102.03"
$ws.Range("E18").Value = "  -4.74%  "
$ws.Range("D19").Value = "'6.75"
$ws.Range("E19").Value = "  -7.31%  "
$ws.Range("D20").Value = "'471.98"
$ws.Range("E20").Value = "  -5.06%  "
$ws.Range("D21").Value = "'14.47"
$ws.Range("E21").Value = "  -5.07%  "
$ws.Range("E22").Value = "  -6.36%  "
$ws.Range("D23").Value = "'7.66"
$ws.Range("E23").Value = "  -4.53%  "
$ws.Range("D24").Value = "'13.41"
$ws.Range("E24").Value = "  -7.78%  "
$ws.Range("D25").Value = "'82.85"
$ws.Range("E25").Value = "  -4.39%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "'2.76"
$ws.Range("E27").Value = "  -8.51%  "
$ws.Range("D28").Value = "'8.34"
$ws.Range("E28").Value = "  -8.03%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.116"
$ws.Range("E29").Value = "  -10.52%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'6.81"
$ws.Range("E30").Value = "  -2.90%  "
$ws.Range("E31").Value = "  -12.37%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "'2.64"
$ws.Range("E33").Value = "  -6.51%  "
$ws.Range("D34").Value = "'26.08"
$ws.Range("E34").Value = "  -5.97%  "
$ws.Range("D35").Value = "'1.11"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("D36").Value = "'5.88"
$ws.Range("E36").Value = "  -7.68%  "
$ws.Range("D37").Value = "'52.41"
$ws.Range("E37").Value = "  -5.60%  "
$ws.Range("D38").Value = "0.0₃0746"
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("D39").Value = "'452.84"
$ws.Range("E39").Value = "  -8.31%  "
$ws.Range("D40").Value = "'2.92"
$ws.Range("E40").Value = "  -15.11%  "
$ws.Range("D41").Value = "'0.0389"
$ws.Range("E41").Value = "  -7.37%  "
$ws.Range("D42").Value = "'0.117"
$ws.Range("E42").Value = "  -9.31%  "
$ws.Range("D43").Value = "'8.28"
$ws.Range("E43").Value = "  -5.29%  "
$ws.Range("D44").Value = "2.829.83"
$ws.Range("E44").Value = "  -5.46%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.26"
$ws.Range("E45").Value = "  -11.76%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.263"
$ws.Range("E46").Value = "  -9.37%  "
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'25.90"
$ws.Range("E49").Value = "  -9.29%  "
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("D51").Value = "'118.61"
$ws.Range("E51").Value = "  -2.06%  "
